# Update "想去人数" (column F) values on the 展览 and 全部类型 sheets
# to reflect the latest generated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row => new F value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2234
$ws1.Range("F3").Value = 101
$ws1.Range("F5").Value = 81
$ws1.Range("F7").Value = 531
$ws1.Range("F8").Value = 495
$ws1.Range("F28").Value = 1044
$ws1.Range("F31").Value = 35
$ws1.Range("F32").Value = 185

# Sheet "全部类型" (row => new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2234
$ws4.Range("F3").Value = 101
$ws4.Range("F5").Value = 81
$ws4.Range("F8").Value = 531
$ws4.Range("F9").Value = 495
$ws4.Range("F29").Value = 1044
$ws4.Range("F32").Value = 35
$ws4.Range("F33").Value = 185
